# Add season-record columns (Wins, Losses, Ties) to the roster sheet.
# Mirrors the commit "Created functions to get season record": for every
# player row, append the team's Wins/Losses/Ties to the right of the
# existing "Unnamed: 28" column (AC), producing new columns AD:AF.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new labels, styled like the existing header cells ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AC1").Copy()
$ws.Range("AD1").PasteSpecial(-4122)   # xlPasteFormats - copy header style (bold/border)

$ws.Range("AE1").Value = "Losses"
$ws.Range("AC1").Copy()
$ws.Range("AE1").PasteSpecial(-4122)

$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AF1").PasteSpecial(-4122)

# --- Data rows (2-45): season record is the same for every player on the roster ---
$wins = 77
$losses = 85
$ties = 0

$lastRow = 45
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins     # column AD
    $ws.Cells.Item($r, 31).Value = $losses   # column AE
    $ws.Cells.Item($r, 32).Value = $ties     # column AF
}
